# Auto-generated Excel COM-interop script to apply market-data price updates
# to the Ragnarok_Profits workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3067.0908
$ws.Range("I62").Value = 2817.5
$ws.Range("K62").Value = 2817.5
$ws.Range("M62").Value = -2193.5

$ws.Range("H65").Value = 3067.0908
$ws.Range("I65").Value = 2817.5
$ws.Range("K65").Value = 14087.5
$ws.Range("M65").Value = -10967.5

$ws.Range("H86").Value = 3236.9333
$ws.Range("J86").Value = 1450
$ws.Range("L86").Value = 1450
$ws.Range("N86").Value = -3696

$ws.Range("H87").Value = 89997.5
$ws.Range("J87").Value = 89997.5
$ws.Range("L87").Value = 89997.5
$ws.Range("N87").Value = -92493.5

$ws.Range("H89").Value = 3236.9333
$ws.Range("J89").Value = 1450
$ws.Range("L89").Value = 7250
$ws.Range("N89").Value = -18482

$ws.Range("H90").Value = 89997.5
$ws.Range("J90").Value = 89997.5
$ws.Range("L90").Value = 269992.5
$ws.Range("N90").Value = -282472.5

$ws.Range("H132").Value = 2214.1875
$ws.Range("I132").Value = 2243.5227
$ws.Range("J132").Value = 1891.5
$ws.Range("K132").Value = 6730.5681
$ws.Range("L132").Value = 5674.5
$ws.Range("M132").Value = -4200.5681
$ws.Range("N132").Value = -10734.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").ClearContents()
$ws.Range("N18").Value = 0

$ws.Range("H32").Value = 13686.158
$ws.Range("I32").Value = 14013.223
$ws.Range("K32").Value = 14013.223
$ws.Range("M32").Value = -13726.223

$ws.Range("H61").Value = 19102170
$ws.Range("I61").Value = 22234874
$ws.Range("K61").Value = 22234874
$ws.Range("M61").Value = -22234662

$ws.Range("H122").Value = 2220.4
$ws.Range("I122").Value = 2182.5715
$ws.Range("K122").Value = 6547.7145
$ws.Range("M122").Value = -4097.7145

$ws.Range("H132").Value = 8339356
$ws.Range("I132").Value = 6973.3335
$ws.Range("K132").Value = 20920.0005
$ws.Range("M132").Value = -18390.0005

$ws.Range("H136").Value = 19102170
$ws.Range("I136").Value = 22234874
$ws.Range("K136").Value = 66704622
$ws.Range("M136").Value = -66702072

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1611.875
$ws.Range("I94").Value = 1858.8
$ws.Range("J94").Value = 1200.3334
$ws.Range("K94").Value = 1858.8
$ws.Range("L94").Value = 1200.3334
$ws.Range("M94").Value = -1407.8
$ws.Range("N94").Value = -2102.3334

$ws.Range("H99").Value = 1836.7693
$ws.Range("I99").Value = 968.2857
$ws.Range("K99").Value = 968.2857
$ws.Range("M99").Value = 529.7143

$ws.Range("I105").Value = 859808.7
$ws.Range("K105").Value = 859808.7
$ws.Range("M105").Value = -858061.7

$ws.Range("H132").Value = 199999
$ws.Range("J132").Value = 199999
$ws.Range("L132").Value = 199999
$ws.Range("N132").Value = -210119

$ws.Range("H134").Value = 10001615
$ws.Range("I134").Value = 1167.8
$ws.Range("J134").Value = 20002062
$ws.Range("K134").Value = 3503.4
$ws.Range("L134").Value = 60006186
$ws.Range("M134").Value = -968.3999999999996
$ws.Range("N134").Value = -60011256

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 76927270
$ws.Range("I31").Value = 90912820
$ws.Range("K31").Value = 90912820
$ws.Range("M31").Value = -90912525

$ws.Range("H32").Value = 333339000
$ws.Range("I32").Value = 333339000
$ws.Range("K32").Value = 333339000
$ws.Range("M32").Value = -333338684

$ws.Range("H34").Value = 76927270
$ws.Range("I34").Value = 90912820
$ws.Range("K34").Value = 90912820
$ws.Range("M34").Value = -90912618

$ws.Range("H132").Value = 2218.5833
$ws.Range("I132").Value = 1840.1428
$ws.Range("J132").Value = 4867.6665
$ws.Range("K132").Value = 5520.428400000001
$ws.Range("L132").Value = 14602.9995
$ws.Range("M132").Value = -2990.428400000001
$ws.Range("N132").Value = -19662.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1699.6666
$ws.Range("I113").Value = 1850.25
$ws.Range("J113").Value = 1624.375
$ws.Range("K113").Value = 5550.75
$ws.Range("L113").Value = 4873.125
$ws.Range("M113").Value = -3380.75
$ws.Range("N113").Value = -9213.125

$ws.Range("H131").Value = 4028.2122
$ws.Range("J131").Value = 3809.7097
$ws.Range("L131").Value = 11429.1291
$ws.Range("N131").Value = -21509.1291

$ws.Range("H134").Value = 8134.381
$ws.Range("I134").Value = 2793.4707
$ws.Range("K134").Value = 8380.4121
$ws.Range("M134").Value = -3310.4121

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 132.61539
$ws.Range("I2").Value = 152.5
$ws.Range("K2").Value = 152.5
$ws.Range("M2").Value = -39.5

$ws.Range("H102").Value = 2733.5789
$ws.Range("I102").Value = 2095.3845
$ws.Range("K102").Value = 2095.3845
$ws.Range("M102").Value = -473.3845000000001

$ws.Range("H126").Value = 17164036
$ws.Range("I126").Value = 23291018
$ws.Range("J126").Value = 8485
$ws.Range("K126").Value = 69873054
$ws.Range("L126").Value = 25455
$ws.Range("M126").Value = -69870584
$ws.Range("N126").Value = -30395

$ws.Range("H132").Value = 5788563.5
$ws.Range("I132").Value = 3612.1924
$ws.Range("J132").Value = 27275526
$ws.Range("K132").Value = 10836.5772
$ws.Range("L132").Value = 81826578
$ws.Range("M132").Value = -8306.5772
$ws.Range("N132").Value = -81831638

$ws.Range("H138").Value = 97499.5
$ws.Range("J138").Value = 97499.5
$ws.Range("L138").Value = 97499.5
$ws.Range("N138").Value = -107779.5

$ws.Range("H139").Value = 134974.75
$ws.Range("J139").Value = 134974.75
$ws.Range("L139").Value = 134974.75
$ws.Range("N139").Value = -145254.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 789.2727
$ws.Range("I22").Value = 675
$ws.Range("J22").Value = 989.25
$ws.Range("K22").Value = 675
$ws.Range("L22").Value = 989.25
$ws.Range("M22").Value = -380
$ws.Range("N22").Value = -1579.25

$ws.Range("H27").Value = 789.2727
$ws.Range("I27").Value = 675
$ws.Range("J27").Value = 989.25
$ws.Range("K27").Value = 675
$ws.Range("L27").Value = 989.25
$ws.Range("M27").Value = -568
$ws.Range("N27").Value = -1203.25

$ws.Range("H46").Value = 886.5
$ws.Range("I46").Value = 705.125
$ws.Range("J46").Value = 1249.25
$ws.Range("K46").Value = 705.125
$ws.Range("L46").Value = 1249.25
$ws.Range("M46").Value = -517.125
$ws.Range("N46").Value = -1625.25

$ws.Range("H61").Value = 4392.923
$ws.Range("I61").Value = 3603.375
$ws.Range("K61").Value = 3603.375
$ws.Range("M61").Value = -3401.375

$ws.Range("H100").Value = 12516301
$ws.Range("I100").Value = 3598.5
$ws.Range("J100").Value = 41712610
$ws.Range("K100").Value = 3598.5
$ws.Range("L100").Value = 41712610
$ws.Range("M100").Value = -3057.5
$ws.Range("N100").Value = -41713692

$ws.Range("H113").Value = 4392.923
$ws.Range("I113").Value = 3603.375
$ws.Range("K113").Value = 3603.375
$ws.Range("M113").Value = -1433.375

$ws.Range("H122").Value = 3382.3538
$ws.Range("I122").Value = 3218.9106
$ws.Range("K122").Value = 9656.731800000001
$ws.Range("M122").Value = -7206.731800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8358.733
$ws.Range("I62").Value = 4337.1113
$ws.Range("K62").Value = 4337.1113
$ws.Range("M62").Value = -3713.1113

$ws.Range("H65").Value = 8358.733
$ws.Range("I65").Value = 4337.1113
$ws.Range("K65").Value = 21685.5565
$ws.Range("M65").Value = -18565.5565

$ws.Range("H107").Value = 2860.1562
$ws.Range("I107").Value = 2112.5454
$ws.Range("J107").Value = 4504.9
$ws.Range("K107").Value = 6337.6362
$ws.Range("L107").Value = 13514.7
$ws.Range("M107").Value = -4417.6362
$ws.Range("N107").Value = -17354.7

$ws.Range("H113").Value = 1351.0769
$ws.Range("I113").Value = 1446.0526
$ws.Range("K113").Value = 4338.1578
$ws.Range("M113").Value = -2168.1578

$ws.Range("H130").Value = 44912.5
$ws.Range("J130").Value = 44912.5
$ws.Range("L130").Value = 44912.5
$ws.Range("N130").Value = -54952.5

$ws.Range("H135").Value = 69784.664
$ws.Range("J135").Value = 69784.664
$ws.Range("L135").Value = 69784.664
$ws.Range("N135").Value = -79924.664

